$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that changed from 45189 (2023-09-20)
# to 45190 (2023-09-21) for every data row (rows 2 through 15).
for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 45189) {
        $cell.Value = 45190
    }
}
